$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.404.32'
$ws.Range('E2').Value = '  -1.21%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.530.68'
$ws.Range('E3').Value = '  +0.18%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.12'
$ws.Range('E5').Value = '  +4.42%  '

# Row 6: 'Solana' -> 'Solana'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.73'
$ws.Range('E6').Value = '  -5.39%  '

# Row 7: 'XRP' -> 'XRP'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.634'
$ws.Range('E7').Value = '  +4.52%  '

# Row 8: 'USDC' -> 'USDC'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.01%  '

# Row 9: 'Cardano' -> 'Cardano'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.639'
$ws.Range('E9').Value = '  +0.73%  '

# Row 10: 'Dogecoin' -> 'Dogecoin'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.158'
$ws.Range('E10').Value = '  +5.14%  '

# Row 11: 'Avalanche' -> 'Avalanche'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '55.70'
$ws.Range('E11').Value = '  +1.27%  '

# Row 12: 'ShibaInu' -> 'ShibaInu'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000279'
$ws.Range('E12').Value = '  +3.53%  '

# Row 13: 'Polkadot' -> 'Polkadot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.31'
$ws.Range('E13').Value = '  -0.78%  '

# Row 14: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.091.01'
$ws.Range('E14').Value = '  +1.05%  '

# Row 15: 'WrappedEther' -> 'WrappedEther'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.525.34'
$ws.Range('E15').Value = '  +0.68%  '

# Row 16: 'TRON' -> 'TRON'
$ws.Range('E16').Value = '  +0.22%  '

# Row 17: 'Chainlink' -> 'Chainlink'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.43'
$ws.Range('E17').Value = '  +1.21%  '

# Row 18: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.350.79'
$ws.Range('E18').Value = '  -1.16%  '

# Row 19: 'Uniswap' -> 'Uniswap'
$ws.Range('E19').Value = '  +1.86%  '

# Row 20: 'Polygon' -> 'Polygon'
$ws.Range('E20').Value = '  +1.41%  '

# Row 21: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range('E21').Value = '  -3.07%  '

# Row 22: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.21'
$ws.Range('E22').Value = '  +7.63%  '

# Row 23: 'Litecoin' -> 'Toncoin'
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.30'
$ws.Range('E23').Value = '  +2.02%  '

# Row 24: 'Toncoin' -> 'Litecoin'
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.81'
$ws.Range('E24').Value = '  +0.83%  '

# Row 25: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range('E25').Value = '  +8.44%  '

# Row 26: 'RenderToken' -> 'RenderToken'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.21'
$ws.Range('E26').Value = '  -0.06%  '

# Row 27: 'ImmutableX' -> 'ImmutableX'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.88'
$ws.Range('E27').Value = '  -1.45%  '

# Row 28: 'Filecoin' -> 'LEO'
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.06'
$ws.Range('E28').Value = '  -1.76%  '

# Row 29: 'EthereumClassic' -> 'Filecoin'
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.11'
$ws.Range('E29').Value = '  +2.52%  '

# Row 30: 'Bittensor' -> 'EthereumClassic'
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.64'
$ws.Range('E30').Value = '  +1.13%  '

# Row 31: 'NEARProtocol' -> 'Bittensor'
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '648.41'
$ws.Range('E31').Value = '  -0.67%  '

# Row 32: 'Cosmos' -> 'NEARProtocol'
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.56'
$ws.Range('E32').Value = '  -1.75%  '

# Row 33: 'Hedera' -> 'Cosmos'
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.78'
$ws.Range('E33').Value = '  +0.61%  '

# Row 34: 'Kaspa' -> 'Hedera'
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.111'
$ws.Range('E34').Value = '  +1.01%  '

# Row 35: 'OKB' -> 'Kaspa'
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.156'
$ws.Range('E35').Value = '  +12.78%  '

# Row 36: 'PEPE' -> 'OKB'
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.79'
$ws.Range('E36').Value = '  +0.21%  '

# Row 37: 'Dai' -> 'PEPE'
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0824'
$ws.Range('E37').Value = '  +0.73%  '

# Row 38: 'InjectiveProtocol' -> 'Dai'
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.27%  '

# Row 39: 'TheGraph' -> 'InjectiveProtocol'
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.60'
$ws.Range('E39').Value = '  -2.95%  '

# Row 40: 'Stacks' -> 'Stacks'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.55'
$ws.Range('E40').Value = '  +7.09%  '

# Row 41: 'Maker' -> 'TheGraph'
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.384'
$ws.Range('E41').Value = '  -1.73%  '

# Row 42: 'FirstDigitalUSD' -> 'Maker'
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.268.84'
$ws.Range('E42').Value = '  +8.35%  '

# Row 43: 'ThetaToken' -> 'FirstDigitalUSD'
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.17%  '

# Row 44: 'ApeXProtocol' -> 'ThetaToken'
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.94'
$ws.Range('E44').Value = '  +1.39%  '

# Row 45: 'Fetch.AI' -> 'Fetch.AI'
$ws.Range('E45').Value = '  -4.20%  '

# Row 46: 'VeChain' -> 'ApeXProtocol'
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.31'
$ws.Range('E46').Value = '  -0.54%  '

# Row 47: 'WEMIXToken' -> 'VeChain'
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0421'
$ws.Range('E47').Value = '  +0.90%  '

# Row 48: 'Stellar' -> 'WEMIXToken'
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.70'
$ws.Range('E48').Value = '  -0.56%  '

# Row 49: 'THORChain' -> 'Stellar'
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.133'
$ws.Range('E49').Value = '  +1.75%  '

# Row 50: 'Monero' -> 'THORChain'
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.72'
$ws.Range('E50').Value = '  -1.15%  '

# Row 51: 'dogwifhat' -> 'Monero'
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '138.66'
$ws.Range('E51').Value = '  -1.49%  '
